$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Range("A1").Value = 'Marca temporal'
$ws.Range("B1").Value = 'Dirección de correo electrónico'
$ws.Range("C1").Value = 'Nombre y apellidos'
$ws.Range("D1").Value = 'Género'
$ws.Range("E1").Value = 'Escuela'
$ws.Range("F1").Value = 'Edad '
$ws.Range("G1").Value = 'Peso (kg sin decimales)'
$ws.Range("H1").Value = 'Tiempo entrenando'
$ws.Range("I1").Value = 'EPS'
$ws.Range("J1").Value = 'Inscripción a '
$ws.Range("K1").Value = 'Categoría formas'
$ws.Range("L1").Value = '¿Es usted cinta negra?'
$ws.Range("M1").Value = 'Telefono Personal'
$ws.Range("N1").Value = 'Nombre completo acudiente'
$ws.Range("O1").Value = 'Teléfono acudiente'

# Row 2
$ws.Range("A2").Value = '6/21/2019 15:49:23'
$ws.Range("B2").Value = 'correo default'
$ws.Range("C2").Value = 'Camila'
$ws.Range("D2").Value = 'Femenino'
$ws.Range("E2").Value = 'Zhang Fei'
$ws.Range("F2").Value = 19
$ws.Range("G2").Value = 50
$ws.Range("H2").Value = '1-3 años'
$ws.Range("I2").Value = 'Comfenalco'
$ws.Range("J2").Value = 'Sanda, Formas'
$ws.Range("K2").Value = 'Forma sin arma'
$ws.Range("L2").Value = 'No'
$ws.Range("M2").Value = 3167527488
$ws.Range("N2").Value = 'Conny'
$ws.Range("O2").Value = 3173694664

# Row 3
$ws.Range("A3").Value = '6/21/2019 15:49:23'
$ws.Range("B3").Value = 'correo default'
$ws.Range("C3").Value = 'Femme 2'
$ws.Range("D3").Value = 'Femenino'
$ws.Range("E3").Value = 'Zhang Fei'
$ws.Range("F3").Value = 20
$ws.Range("G3").Value = 60
$ws.Range("H3").Value = '1-3 años'
$ws.Range("I3").Value = 'Eps'
$ws.Range("J3").Value = 'Sanda'
$ws.Range("L3").Value = 'No'
$ws.Range("M3").Value = 3167527488
$ws.Range("N3").Value = 'Acudientes'
$ws.Range("O3").Value = 3173694664

# Row 4
$ws.Range("A4").Value = '6/23/2019 15:49:23'
$ws.Range("B4").Value = 'correo default'
$ws.Range("C4").Value = 'Femme 3'
$ws.Range("D4").Value = 'Femenino'
$ws.Range("E4").Value = 'Zhang Fei'
$ws.Range("F4").Value = 21
$ws.Range("G4").Value = 65
$ws.Range("H4").Value = '1-3 años'
$ws.Range("I4").Value = 'Eps'
$ws.Range("J4").Value = 'Sanda'
$ws.Range("L4").Value = 'No'
$ws.Range("M4").Value = 3167527488
$ws.Range("N4").Value = 'Acudientes'
$ws.Range("O4").Value = 3173694666

# Row 5
$ws.Range("A5").Value = '6/27/2019 15:49:23'
$ws.Range("B5").Value = 'correo default'
$ws.Range("C5").Value = 'Femme 4'
$ws.Range("D5").Value = 'Femenino'
$ws.Range("E5").Value = 'Zhang Fei'
$ws.Range("F5").Value = 22
$ws.Range("G5").Value = 70
$ws.Range("H5").Value = '0-1 año'
$ws.Range("I5").Value = 'Eps'
$ws.Range("J5").Value = 'Sanda'
$ws.Range("L5").Value = 'No'
$ws.Range("M5").Value = 3167527488
$ws.Range("N5").Value = 'Acudientes'
$ws.Range("O5").Value = 3173694670

# Row 6
$ws.Range("A6").Value = '6/21/2019 15:49:23'
$ws.Range("B6").Value = 'correo default'
$ws.Range("C6").Value = 'Femme 5'
$ws.Range("D6").Value = 'Femenino'
$ws.Range("E6").Value = 'Zhang Fei'
$ws.Range("F6").Value = 23
$ws.Range("G6").Value = 67
$ws.Range("H6").Value = '0-1 año'
$ws.Range("I6").Value = 'Eps'
$ws.Range("J6").Value = 'Sanda'
$ws.Range("L6").Value = 'No'
$ws.Range("M6").Value = 3167527488
$ws.Range("N6").Value = 'Acudientes'
$ws.Range("O6").Value = 3173694671

# Row 7
$ws.Range("A7").Value = '6/21/2019 15:49:23'
$ws.Range("B7").Value = 'correo default'
$ws.Range("C7").Value = 'Femme 6'
$ws.Range("D7").Value = 'Femenino'
$ws.Range("E7").Value = 'Zhang Fei'
$ws.Range("F7").Value = 24
$ws.Range("G7").Value = 68
$ws.Range("H7").Value = '1-3 años'
$ws.Range("I7").Value = 'Eps'
$ws.Range("J7").Value = 'Sanda'
$ws.Range("L7").Value = 'No'
$ws.Range("M7").Value = 3167527488
$ws.Range("N7").Value = 'Acudientes'
$ws.Range("O7").Value = 3173694673

# Row 8
$ws.Range("A8").Value = '6/23/2019 15:49:24'
$ws.Range("B8").Value = 'correo default'
$ws.Range("C8").Value = 'Femme 7'
$ws.Range("D8").Value = 'Femenino'
$ws.Range("E8").Value = 'Zhang Fei'
$ws.Range("F8").Value = 25
$ws.Range("G8").Value = 50
$ws.Range("H8").Value = '1-3 años'
$ws.Range("I8").Value = 'Eps'
$ws.Range("J8").Value = 'Sanda'
$ws.Range("L8").Value = 'No'
$ws.Range("M8").Value = 3167527488
$ws.Range("N8").Value = 'Acudientes'
$ws.Range("O8").Value = 3173694675

# Row 9
$ws.Range("A9").Value = '6/27/2019 15:49:24'
$ws.Range("B9").Value = 'correo default'
$ws.Range("C9").Value = 'Femme 8'
$ws.Range("D9").Value = 'Femenino'
$ws.Range("E9").Value = 'Zhang Fei'
$ws.Range("F9").Value = 26
$ws.Range("G9").Value = 69
$ws.Range("H9").Value = '0-1 año'
$ws.Range("I9").Value = 'Eps'
$ws.Range("J9").Value = 'Sanda'
$ws.Range("L9").Value = 'No'
$ws.Range("M9").Value = 3167527488
$ws.Range("N9").Value = 'Acudientes'
$ws.Range("O9").Value = 3173694677

# Row 10
$ws.Range("A10").Value = '6/21/2019 15:49:23'
$ws.Range("B10").Value = 'correo default'
$ws.Range("C10").Value = 'Femme 9'
$ws.Range("D10").Value = 'Femenino'
$ws.Range("E10").Value = 'Zhang Fei'
$ws.Range("F10").Value = 27
$ws.Range("G10").Value = 80
$ws.Range("H10").Value = '0-1 año'
$ws.Range("I10").Value = 'Eps'
$ws.Range("J10").Value = 'Sanda'
$ws.Range("L10").Value = 'No'
$ws.Range("M10").Value = 3167527488
$ws.Range("N10").Value = 'Acudientes'
$ws.Range("O10").Value = 3173694679

# Row 11
$ws.Range("A11").Value = '6/21/2019 15:49:23'
$ws.Range("B11").Value = 'correo default'
$ws.Range("C11").Value = 'Femme 10'
$ws.Range("D11").Value = 'Femenino'
$ws.Range("E11").Value = 'Zhang Fei'
$ws.Range("F11").Value = 28
$ws.Range("G11").Value = 70
$ws.Range("H11").Value = '1-3 años'
$ws.Range("I11").Value = 'Eps'
$ws.Range("J11").Value = 'Sanda'
$ws.Range("L11").Value = 'No'
$ws.Range("M11").Value = 3167527488
$ws.Range("N11").Value = 'Acudientes'
$ws.Range("O11").Value = 3173694681

# Row 12
$ws.Range("A12").Value = '6/23/2019 15:49:25'
$ws.Range("B12").Value = 'correo default'
$ws.Range("C12").Value = 'Femme 11'
$ws.Range("D12").Value = 'Femenino'
$ws.Range("E12").Value = 'Zhang Fei'
$ws.Range("F12").Value = 29
$ws.Range("G12").Value = 61
$ws.Range("H12").Value = '1-3 años'
$ws.Range("I12").Value = 'Eps'
$ws.Range("J12").Value = 'Sanda'
$ws.Range("L12").Value = 'No'
$ws.Range("M12").Value = 3167527488
$ws.Range("N12").Value = 'Acudientes'
$ws.Range("O12").Value = 3173694683

# Row 13
$ws.Range("A13").Value = '6/27/2019 15:49:25'
$ws.Range("B13").Value = 'correo default'
$ws.Range("C13").Value = 'Femme 12'
$ws.Range("D13").Value = 'Femenino'
$ws.Range("E13").Value = 'Zhang Fei'
$ws.Range("F13").Value = 30
$ws.Range("G13").Value = 71
$ws.Range("H13").Value = '0-1 año'
$ws.Range("I13").Value = 'Eps'
$ws.Range("J13").Value = 'Sanda'
$ws.Range("L13").Value = 'No'
$ws.Range("M13").Value = 3167527488
$ws.Range("N13").Value = 'Acudientes'
$ws.Range("O13").Value = 3173694685

# Row 14
$ws.Range("A14").Value = '6/21/2019 15:49:23'
$ws.Range("B14").Value = 'correo default'
$ws.Range("C14").Value = 'Femme 13'
$ws.Range("D14").Value = 'Femenino'
$ws.Range("E14").Value = 'Zhang Fei'
$ws.Range("F14").Value = 31
$ws.Range("G14").Value = 30
$ws.Range("H14").Value = '0-1 año'
$ws.Range("I14").Value = 'Eps'
$ws.Range("J14").Value = 'Sanda'
$ws.Range("L14").Value = 'No'
$ws.Range("M14").Value = 3167527488
$ws.Range("N14").Value = 'Acudientes'
$ws.Range("O14").Value = 3173694687

# Row 15
$ws.Range("A15").Value = '6/21/2019 15:49:23'
$ws.Range("B15").Value = 'correo default'
$ws.Range("C15").Value = 'Femme 14'
$ws.Range("D15").Value = 'Femenino'
$ws.Range("E15").Value = 'Zhang Fei'
$ws.Range("F15").Value = 32
$ws.Range("G15").Value = 72
$ws.Range("H15").Value = '1-3 años'
$ws.Range("I15").Value = 'Eps'
$ws.Range("J15").Value = 'Sanda'
$ws.Range("L15").Value = 'No'
$ws.Range("M15").Value = 3167527488
$ws.Range("N15").Value = 'Acudientes'
$ws.Range("O15").Value = 3173694689

# Row 16
$ws.Range("A16").Value = '6/23/2019 15:49:26'
$ws.Range("B16").Value = 'correo default'
$ws.Range("C16").Value = 'Femme 15'
$ws.Range("D16").Value = 'Femenino'
$ws.Range("E16").Value = 'Zhang Fei'
$ws.Range("F16").Value = 33
$ws.Range("G16").Value = 85
$ws.Range("H16").Value = '1-3 años'
$ws.Range("I16").Value = 'Eps'
$ws.Range("J16").Value = 'Formas'
$ws.Range("K16").Value = 'Forma sin arma, Forma con arma'
$ws.Range("L16").Value = 'No'
$ws.Range("M16").Value = 3167527488
$ws.Range("N16").Value = 'Acudientes'
$ws.Range("O16").Value = 3173694691

# Row 17
$ws.Range("A17").Value = '6/27/2019 15:49:26'
$ws.Range("B17").Value = 'correo default'
$ws.Range("C17").Value = 'Femme 16'
$ws.Range("D17").Value = 'Femenino'
$ws.Range("E17").Value = 'Zhang Fei'
$ws.Range("F17").Value = 34
$ws.Range("G17").Value = 73
$ws.Range("H17").Value = '0-1 año'
$ws.Range("I17").Value = 'Eps'
$ws.Range("J17").Value = 'Formas'
$ws.Range("K17").Value = 'Forma sin arma'
$ws.Range("L17").Value = 'No'
$ws.Range("M17").Value = 3167527488
$ws.Range("N17").Value = 'Acudientes'
$ws.Range("O17").Value = 3173694693

# Row 18
$ws.Range("A18").Value = '6/21/2019 15:49:23'
$ws.Range("B18").Value = 'correo default'
$ws.Range("C18").Value = 'Femme 17'
$ws.Range("D18").Value = 'Femenino'
$ws.Range("E18").Value = 'Zhang Fei'
$ws.Range("F18").Value = 35
$ws.Range("G18").Value = 40
$ws.Range("H18").Value = '0-1 año'
$ws.Range("I18").Value = 'Eps'
$ws.Range("J18").Value = 'Formas'
$ws.Range("K18").Value = 'Forma con arma'
$ws.Range("L18").Value = 'No'
$ws.Range("M18").Value = 3167527488
$ws.Range("N18").Value = 'Acudientes'
$ws.Range("O18").Value = 3173694695

# Row 19
$ws.Range("A19").Value = '6/21/2019 15:49:23'
$ws.Range("B19").Value = 'correo default'
$ws.Range("C19").Value = 'Femme 18'
$ws.Range("D19").Value = 'Femenino'
$ws.Range("E19").Value = 'Zhang Fei'
$ws.Range("F19").Value = 36
$ws.Range("G19").Value = 74
$ws.Range("H19").Value = '1-3 años'
$ws.Range("I19").Value = 'Eps'
$ws.Range("J19").Value = 'Formas'
$ws.Range("K19").Value = 'Forma con arma'
$ws.Range("L19").Value = 'No'
$ws.Range("M19").Value = 3167527488
$ws.Range("N19").Value = 'Acudientes'
$ws.Range("O19").Value = 3173694697

# Row 20
$ws.Range("A20").Value = '6/21/2019 15:49:23'
$ws.Range("B20").Value = 'correo default'
$ws.Range("C20").Value = 'Femme 19'
$ws.Range("D20").Value = 'Femenino'
$ws.Range("E20").Value = 'Zhang Fei'
$ws.Range("F20").Value = 37
$ws.Range("G20").Value = 44
$ws.Range("H20").Value = '0-1 año'
$ws.Range("I20").Value = 'Eps'
$ws.Range("J20").Value = 'Formas'
$ws.Range("K20").Value = 'Forma con arma'
$ws.Range("L20").Value = 'No'
$ws.Range("M20").Value = 3167527488
$ws.Range("N20").Value = 'Acudientes'
$ws.Range("O20").Value = 3173694699

# Column width adjustments (approximate best-fit; engine quantizes to 1/6 char units)
$ws.Columns.Item(9).ColumnWidth = 10.59
$ws.Columns.Item(11).ColumnWidth = 11.42
$ws.Columns.Item(12).ColumnWidth = 19.59
$ws.Columns.Item(13).ColumnWidth = 16.42
$ws.Columns.Item(14).ColumnWidth = 11.59
$ws.Columns.Item(15).ColumnWidth = 10.59

# Update the _xlnm._FilterDatabase defined name to extend through column O
foreach ($n in $wb.Names) {
    if ($n.Name -like "*FilterDatabase*") {
        $n.RefersTo = "=RegistroTest3!`$A`$1:`$O`$7"
    }
}

# Set the active selection to I15 (matches post-edit cursor position)
$ws.Range("I15").Select()
